$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 was stored as an inline/shared text string "3"; the new version stores
# it as a genuine numeric value 3.
$ws.Range("B3").Value = 3

# Add a new row 4, mostly mirroring row 3 but for a different annotation
# (politeness_score 2, a shorter polite_expressions excerpt, and
# sentence_purpose "APC"). Column B on this row must stay textual ("2"),
# matching the source data's inline-string type.
$ws.Range("A4").Value = "Ying Tang"

$b4 = $ws.Cells.Item(4, 2)
$b4.NumberFormat = "@"
$b4.Value = "2"
$b4.Style = "Normal"

$ws.Range("C4").Value = " should be mentioned"
$ws.Range("D4").Value = "APC"
$ws.Range("E4").Value = "THE"
$ws.Range("F4").Value = "c8048836-24fe-4e27-95aa-c7cfb58ac155"
$ws.Range("G4").Value = "rkc_hGb0Z_annotated.xlsx"
$ws.Range("H4").Value = "The structure of the global policies used in the experiments should be mentioned somewhere."
